# Atualizada planilha de calibração do pêndulo
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet2: rework the calibration table (rows 5-15) ---
# Old layout had a "T" / "s" header in row 5, raw period samples in column B
# (rows 6-15) and instrument metadata scattered in columns D/E every other
# row. New layout keeps only the instrument metadata, compacted into
# columns B/C starting at row 5, and drops the raw sample rows entirely.
$ws2.Range("A5:E15").ClearContents()

$ws2.Range("B5").Value = "regua"
$ws2.Range("C5").Value = "30cm"

$ws2.Range("B6").Value = "res"
$ws2.Range("C6").Value = "0.1cm"

$ws2.Range("B8").Value = "cronometro"
$ws2.Range("C8").Value = "iphone"

$ws2.Range("B10").Value = "multimetro"
$ws2.Range("C10").Value = "Fluke 87"

# Sheet1 keeps its own last selection (unchanged) but is no longer the
# tab shown when the workbook is reopened.
$ws1.Activate() | Out-Null
$ws1.Range("E17").Select() | Out-Null

# --- Window / view state: Sheet2 becomes the active, selected tab ---
$ws2.Activate() | Out-Null
$ws2.Range("E7").Select() | Out-Null

# Resize the saved workbook window (best effort; cosmetic UI metric).
$win = $excel.ActiveWindow
$win.Width = 28800
$win.Height = 16240
